# Atualização de bases das ligas, do dia: 28-04-2024 às 15:37
# Applies the odds-database refresh described by the commit diff:
#  - Row 36 and Row 37 (played matches) swap their full data (everything
#    except the running index in column A).
#  - Row 171 and Row 172 (upcoming fixtures) get updated kickoff dates,
#    matchups and odds, and their "id" (column B) shared-string text is
#    refreshed while staying text-typed (not coerced to a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData($Sheet, $Row, $Values) {
    foreach ($col in $Values.Keys) {
        $Sheet.Cells.Item($Row, $col).Value = $Values[$col]
    }
}

# --- Row 36 (id col A = 34) becomes what row 37 used to hold ------------
Set-RowData $ws 36 @{
    2  = 6864629                # B  id
    5  = "Borac Banja Luka"     # E  HomeTeam
    6  = "NK Posusje"           # F  AwayTeam
    7  = 1                      # G  FTHG
    8  = 0                      # H  FTAG
    10 = 1.363                  # J  oddH_op
    11 = 4.5                    # K  oddD_op
    12 = 6.5                    # L  oddA_op
    13 = 1.363                  # M  oddH
    14 = 4.2                    # N  oddD
    15 = 6.5                    # O  oddA
    17 = 1.95                   # Q  oddAHH
    18 = 1.85                   # R  oddAHA
    19 = 2.5                    # S  AhOU
    20 = 1.925                  # T  oddAHOver
    21 = 1.875                  # U  oddAHUnder
    22 = 0.363                  # V  PLH
    26 = 0.425                  # Z  PL_Ahh
    27 = -1                     # AA PL_Aha
    28 = 0.875                  # AB PL_AhOver
}

# --- Row 37 (id col A = 35) becomes what row 36 used to hold ------------
Set-RowData $ws 37 @{
    2  = 6865299
    5  = "Siroki Brijeg"
    6  = "Zvijezda 09"
    7  = 2
    8  = 1
    10 = 1.25
    11 = 5.5
    12 = 8
    13 = 1.4
    14 = 4.75
    15 = 5.75
    17 = 1.9
    18 = 1.9
    19 = 2.75
    20 = 1.85
    21 = 1.95
    22 = 0.3999999999999999
    26 = 0.45
    27 = 0.425
    28 = -0.5
}

# --- Row 171: re-scheduled fixture --------------------------------------
# Force the numeric-looking id to stay text (matches shared-string type in
# the source file) without leaving the cell's number format altered.
$ws.Cells.Item(171, 2).NumberFormat = "@"
$ws.Cells.Item(171, 2).Value = "7952461"
$ws.Cells.Item(171, 2).Style = "Normal"

Set-RowData $ws 171 @{
    4  = 45410.5                # D  Date (serial, keeps existing date format)
    5  = "NK Posusje"           # E  HomeTeam
    6  = "Zrinjski Mostar"      # F  AwayTeam
    10 = 3.5                    # J  oddH_op
    11 = 3.6                    # K  oddD_op
    12 = 1.8                    # L  oddA_op
    13 = 4.5                    # M  oddH
    15 = 1.666                  # O  oddA
    16 = 0.75                   # P  Ah
    17 = 1.85                   # Q  oddAHH
    18 = 1.95                   # R  oddAHA
    20 = 1.975                  # T  oddAHOver
    21 = 1.825                  # U  oddAHUnder
}

# --- Row 172: re-scheduled fixture --------------------------------------
$ws.Cells.Item(172, 2).NumberFormat = "@"
$ws.Cells.Item(172, 2).Value = "7952764"
$ws.Cells.Item(172, 2).Style = "Normal"

Set-RowData $ws 172 @{
    4  = 45410.65625
    5  = "FK Sarajevo"
    6  = "Velez Mostar"
    10 = 1.533
    11 = 3.5
    12 = 5.5
    13 = 1.727
    14 = 3.5
    15 = 4
    16 = -0.75
    17 = 2
    18 = 1.8
    19 = 2.25
    20 = 1.825
    21 = 1.975
}
